# Update res_line/loading_percent values for the 380 kV case (Case_5_208)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 16.42638983778676
$ws.Cells.Item(2, 3).Value = 10.29331459755666
$ws.Cells.Item(2, 4).Value = 7.30167665320554
$ws.Cells.Item(2, 5).Value = 16.45433944982152
$ws.Cells.Item(2, 6).Value = 43.42741042073344
$ws.Cells.Item(3, 2).Value = 15.82553249713264
$ws.Cells.Item(3, 3).Value = 9.688831710892437
$ws.Cells.Item(3, 4).Value = 7.136393549518607
$ws.Cells.Item(3, 5).Value = 15.51418462611142
$ws.Cells.Item(3, 6).Value = 41.92871716270902
$ws.Cells.Item(4, 2).Value = 15.45323296866087
$ws.Cells.Item(4, 3).Value = 9.302815878495984
$ws.Cells.Item(4, 4).Value = 7.033694688316622
$ws.Cells.Item(4, 5).Value = 14.91370740114557
$ws.Cells.Item(4, 6).Value = 40.99073186571778
$ws.Cells.Item(5, 2).Value = 15.30097127547397
$ws.Cells.Item(5, 3).Value = 9.141972636156947
$ws.Cells.Item(5, 4).Value = 6.991576212751836
$ws.Cells.Item(5, 5).Value = 14.66344036576686
$ws.Cells.Item(5, 6).Value = 40.60455364992115
$ws.Cells.Item(6, 2).Value = 15.27566457246559
$ws.Cells.Item(6, 3).Value = 9.115057593640415
$ws.Cells.Item(6, 4).Value = 6.984567413363458
$ws.Cells.Item(6, 5).Value = 14.62155641059596
$ws.Cells.Item(6, 6).Value = 40.54020732386432
$ws.Cells.Item(7, 2).Value = 15.4511813125667
$ws.Cells.Item(7, 3).Value = 9.300660726981791
$ws.Cells.Item(7, 4).Value = 7.033127697937479
$ws.Cells.Item(7, 5).Value = 14.91035435876094
$ws.Cells.Item(7, 6).Value = 40.98553895853991
$ws.Cells.Item(8, 2).Value = 16.22008705086596
$ws.Cells.Item(8, 3).Value = 10.08808041295122
$ws.Cells.Item(8, 4).Value = 7.244961183216043
$ws.Cells.Item(8, 5).Value = 16.13513958936087
$ws.Cells.Item(8, 6).Value = 42.91465571275551
$ws.Cells.Item(9, 2).Value = 17.68917423310451
$ws.Cells.Item(9, 3).Value = 11.65565728110506
$ws.Cells.Item(9, 4).Value = 7.649156145900511
$ws.Cells.Item(9, 5).Value = 18.42300014397892
$ws.Cells.Item(9, 6).Value = 46.53555840643371
$ws.Cells.Item(10, 2).Value = 18.73062856694719
$ws.Cells.Item(10, 3).Value = 12.72416089805899
$ws.Cells.Item(10, 4).Value = 7.937297104013416
$ws.Cells.Item(10, 5).Value = 20.09240612184273
$ws.Cells.Item(10, 6).Value = 49.07208829654715
$ws.Cells.Item(11, 2).Value = 19.19367544082326
$ws.Cells.Item(11, 3).Value = 13.18385019090909
$ws.Cells.Item(11, 4).Value = 8.066081172719578
$ws.Cells.Item(11, 5).Value = 20.81095225300755
$ws.Cells.Item(11, 6).Value = 50.19475926645455
$ws.Cells.Item(12, 2).Value = 19.36730150894325
$ws.Cells.Item(12, 3).Value = 13.3541507046749
$ws.Cells.Item(12, 4).Value = 8.114489260347284
$ws.Cells.Item(12, 5).Value = 21.07721539202315
$ws.Cells.Item(12, 6).Value = 50.61508653765809
$ws.Cells.Item(13, 2).Value = 19.32998712384931
$ws.Cells.Item(13, 3).Value = 13.31764096272882
$ws.Cells.Item(13, 4).Value = 8.104080161986687
$ws.Cells.Item(13, 5).Value = 21.02012967387631
$ws.Cells.Item(13, 6).Value = 50.52477981594377
$ws.Cells.Item(14, 2).Value = 19.20799517259175
$ws.Cells.Item(14, 3).Value = 13.19793643458213
$ws.Cells.Item(14, 4).Value = 8.07007110694499
$ws.Cells.Item(14, 5).Value = 20.83297457037248
$ws.Cells.Item(14, 6).Value = 50.22943779388444
$ws.Cells.Item(15, 2).Value = 19.13304269629683
$ws.Cells.Item(15, 3).Value = 13.12412308575819
$ws.Cells.Item(15, 4).Value = 8.049191890136083
$ws.Cells.Item(15, 5).Value = 20.7175780972205
$ws.Cells.Item(15, 6).Value = 50.04789795263731
$ws.Cells.Item(16, 2).Value = 18.70013593250552
$ws.Cells.Item(16, 3).Value = 12.69358994218468
$ws.Cells.Item(16, 4).Value = 7.928832017558936
$ws.Cells.Item(16, 5).Value = 20.04462861492908
$ws.Cells.Item(16, 6).Value = 48.99806319788954
$ws.Cells.Item(17, 2).Value = 18.43168260467048
$ws.Cells.Item(17, 3).Value = 12.42272897323808
$ws.Cells.Item(17, 4).Value = 7.854385926304401
$ws.Cells.Item(17, 5).Value = 19.62135747802872
$ws.Cells.Item(17, 6).Value = 48.34579288210324
$ws.Cells.Item(18, 2).Value = 18.27627591557962
$ws.Cells.Item(18, 3).Value = 12.26445279453743
$ws.Cells.Item(18, 4).Value = 7.811352361224981
$ws.Cells.Item(18, 5).Value = 19.37405318474894
$ws.Cells.Item(18, 6).Value = 47.96770302184997
$ws.Cells.Item(19, 2).Value = 18.22349204673524
$ws.Cells.Item(19, 3).Value = 12.21043598561017
$ws.Cells.Item(19, 4).Value = 7.796746104326484
$ws.Cells.Item(19, 5).Value = 19.28965765983189
$ws.Cells.Item(19, 6).Value = 47.83919683658041
$ws.Cells.Item(20, 2).Value = 18.46036472424897
$ws.Cells.Item(20, 3).Value = 12.45181954634602
$ws.Cells.Item(20, 4).Value = 7.862333226540804
$ws.Cells.Item(20, 5).Value = 19.6668135541616
$ws.Cells.Item(20, 6).Value = 48.41553291880169
$ws.Cells.Item(21, 2).Value = 19.24387515290404
$ws.Cells.Item(21, 3).Value = 13.23319881649802
$ws.Cells.Item(21, 4).Value = 8.080070379528696
$ws.Cells.Item(21, 5).Value = 20.88810456907759
$ws.Cells.Item(21, 6).Value = 50.3163196081298
$ws.Cells.Item(22, 2).Value = 19.74585311027923
$ws.Cells.Item(22, 3).Value = 13.72188835381062
$ws.Cells.Item(22, 4).Value = 8.220265545473767
$ws.Cells.Item(22, 5).Value = 21.65229987435612
$ws.Cells.Item(22, 6).Value = 51.53047303727401
$ws.Cells.Item(23, 2).Value = 19.47891475148388
$ws.Cells.Item(23, 3).Value = 13.46307063694064
$ws.Cells.Item(23, 4).Value = 8.145643027621672
$ws.Cells.Item(23, 5).Value = 21.24753019877283
$ws.Cells.Item(23, 6).Value = 50.8851241141432
$ws.Cells.Item(24, 2).Value = 18.4474008547955
$ws.Cells.Item(24, 3).Value = 12.43867565670082
$ws.Cells.Item(24, 4).Value = 7.858740977155679
$ws.Cells.Item(24, 5).Value = 19.64627520016255
$ws.Cells.Item(24, 6).Value = 48.38401305149029
$ws.Cells.Item(25, 2).Value = 17.29746576781898
$ws.Cells.Item(25, 3).Value = 11.23866050126309
$ws.Cells.Item(25, 4).Value = 7.541222888398656
$ws.Cells.Item(25, 5).Value = 17.77169398848133
$ws.Cells.Item(25, 6).Value = 45.57610727184977
